$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4220.3335
$ws.Range("I2").Value = 1483.1666
$ws.Range("K2").Value = 1483.1666
$ws.Range("M2").Value = -1370.1666
$ws.Range("H17").Value = 3011536.5
$ws.Range("J17").Value = 3011536.5
$ws.Range("L17").Value = 9034609.5
$ws.Range("N17").Value = -9034945.5
$ws.Range("H40").Value = 2789.125
$ws.Range("J40").Value = 2883.818
$ws.Range("L40").Value = 2883.818
$ws.Range("N40").Value = -3233.818
$ws.Range("H86").Value = 58827148
$ws.Range("I86").Value = 111114150
$ws.Range("K86").Value = 111114150
$ws.Range("M86").Value = -111113027
$ws.Range("H88").Value = 4553.6924
$ws.Range("I88").Value = 2697.7
$ws.Range("K88").Value = 2697.7
$ws.Range("M88").Value = -2291.7
$ws.Range("H89").Value = 58827148
$ws.Range("I89").Value = 111114150
$ws.Range("K89").Value = 555570750
$ws.Range("M89").Value = -555565134
$ws.Range("H91").Value = 4553.6924
$ws.Range("I91").Value = 2697.7
$ws.Range("K91").Value = 2697.7
$ws.Range("M91").Value = -1293.7
$ws.Range("H92").Value = 1420904.1
$ws.Range("I92").Value = 601301.75
$ws.Range("J92").Value = 4465141.5
$ws.Range("K92").Value = 601301.75
$ws.Range("L92").Value = 4465141.5
$ws.Range("M92").Value = -600053.75
$ws.Range("N92").Value = -4467637.5
$ws.Range("H97").Value = 5444.143
$ws.Range("J97").Value = 5444.143
$ws.Range("L97").Value = 16332.429
$ws.Range("N97").Value = -17324.429
$ws.Range("H116").Value = 2999.6667
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 3249.5
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 3249.5
$ws.Range("N116").Value = -10133.5
$ws.Range("H126").Value = 136663
$ws.Range("J126").Value = 136663
$ws.Range("L126").Value = 136663
$ws.Range("N126").Value = -146543
$ws.Range("H127").Value = 638.8182
$ws.Range("I127").Value = 638.8182
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1916.4546
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3043.5454
$ws.Range("H132").Value = 2528.5833
$ws.Range("I132").Value = 2564.8462
$ws.Range("J132").Value = 2485.7273
$ws.Range("K132").Value = 7694.5386
$ws.Range("L132").Value = 7457.1819
$ws.Range("M132").Value = -5164.5386
$ws.Range("N132").Value = -12517.1819
$ws.Range("H137").Value = 1392954
$ws.Range("I137").Value = 3302.75
$ws.Range("J137").Value = 3576691.5
$ws.Range("K137").Value = 9908.25
$ws.Range("L137").Value = 10730074.5
$ws.Range("M137").Value = -7358.25
$ws.Range("N137").Value = -10735174.5
$ws.Range("H138").Value = 4695
$ws.Range("I138").Value = 5454.048
$ws.Range("K138").Value = 16362.144
$ws.Range("M138").Value = -11222.144
$ws.Range("M116").Value = 942
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2614.7932
$ws.Range("I32").Value = 984.0137
$ws.Range("K32").Value = 984.0137
$ws.Range("M32").Value = -697.0137
$ws.Range("H45").Value = 39252.965
$ws.Range("I45").Value = 47947.09
$ws.Range("K45").Value = 47947.09
$ws.Range("M45").Value = -47570.09
$ws.Range("H61").Value = 1224159.8
$ws.Range("I61").Value = 27727.85
$ws.Range("K61").Value = 27727.85
$ws.Range("M61").Value = -27515.85
$ws.Range("H74").Value = 560301.3
$ws.Range("J74").Value = 1314940.8
$ws.Range("L74").Value = 1314940.8
$ws.Range("N74").Value = -1316688.8
$ws.Range("H77").Value = 560301.3
$ws.Range("J77").Value = 1314940.8
$ws.Range("L77").Value = 6574704
$ws.Range("N77").Value = -6583440
$ws.Range("H110").Value = 1160.2
$ws.Range("I110").Value = 1160.2
$ws.Range("K110").Value = 1160.2
$ws.Range("M110").Value = 884.8
$ws.Range("H132").Value = 2771.9697
$ws.Range("I132").Value = 2289
$ws.Range("K132").Value = 6867
$ws.Range("M132").Value = -4337
$ws.Range("H136").Value = 1224159.8
$ws.Range("I136").Value = 27727.85
$ws.Range("K136").Value = 83183.54999999999
$ws.Range("M136").Value = -80633.54999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1581.0938
$ws.Range("I20").Value = 1378
$ws.Range("J20").Value = 2100.111
$ws.Range("K20").Value = 1378
$ws.Range("L20").Value = 2100.111
$ws.Range("M20").Value = -1131
$ws.Range("N20").Value = -2594.111
$ws.Range("H99").Value = 11612.333
$ws.Range("I99").Value = 10504.177
$ws.Range("K99").Value = 10504.177
$ws.Range("M99").Value = -9006.177
$ws.Range("H134").Value = 30002768
$ws.Range("I134").Value = 2599.524
$ws.Range("J134").Value = 100003160
$ws.Range("K134").Value = 7798.572
$ws.Range("L134").Value = 300009480
$ws.Range("M134").Value = -5263.572
$ws.Range("N134").Value = -300014550
$ws.Range("H135").Value = 149993
$ws.Range("J135").Value = 149993
$ws.Range("L135").Value = 149993
$ws.Range("N135").Value = -160133

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2121.3584
$ws.Range("I31").Value = 2072.4666
$ws.Range("J31").Value = 2140.658
$ws.Range("K31").Value = 2072.4666
$ws.Range("L31").Value = 2140.658
$ws.Range("M31").Value = -1777.4666
$ws.Range("N31").Value = -2730.658
$ws.Range("H34").Value = 2121.3584
$ws.Range("I34").Value = 2072.4666
$ws.Range("J34").Value = 2140.658
$ws.Range("K34").Value = 2072.4666
$ws.Range("L34").Value = 2140.658
$ws.Range("M34").Value = -1870.4666
$ws.Range("N34").Value = -2544.658
$ws.Range("H62").Value = 7121.4546
$ws.Range("I62").Value = 5433.6
$ws.Range("K62").Value = 5433.6
$ws.Range("M62").Value = -4809.6
$ws.Range("H65").Value = 7121.4546
$ws.Range("I65").Value = 5433.6
$ws.Range("K65").Value = 27168
$ws.Range("M65").Value = -24048
$ws.Range("H122").Value = 1951.6666
$ws.Range("I122").Value = 1899.2
$ws.Range("K122").Value = 5697.6
$ws.Range("M122").Value = -3247.6
$ws.Range("H132").Value = 11113653
$ws.Range("I132").Value = 1786.15
$ws.Range("K132").Value = 5358.450000000001
$ws.Range("M132").Value = -2828.450000000001
$ws.Range("H141").Value = 193135.08
$ws.Range("J141").Value = 236036.22
$ws.Range("L141").Value = 236036.22
$ws.Range("N141").Value = -246396.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11991.066
$ws.Range("I3").Value = 6765.222
$ws.Range("J3").Value = 19829.834
$ws.Range("K3").Value = 20295.666
$ws.Range("L3").Value = 59489.50199999999
$ws.Range("M3").Value = -20183.666
$ws.Range("N3").Value = -59713.50199999999
$ws.Range("H56").Value = 10996409
$ws.Range("I56").Value = 10996409
$ws.Range("K56").Value = 10996409
$ws.Range("M56").Value = -10995879
$ws.Range("J81").Value = 10791
$ws.Range("L81").Value = 32373
$ws.Range("N81").Value = -34619
$ws.Range("H82").Value = 18999
$ws.Range("J82").Value = 18999
$ws.Range("L82").Value = 56997
$ws.Range("J84").Value = 10791
$ws.Range("L84").Value = 97119
$ws.Range("N84").Value = -108351
$ws.Range("H85").Value = 18999
$ws.Range("J85").Value = 18999
$ws.Range("L85").Value = 56997
$ws.Range("H120").Value = 22428
$ws.Range("I120").Value = 22349.5
$ws.Range("J120").Value = 22459.4
$ws.Range("K120").Value = 67048.5
$ws.Range("L120").Value = 67378.20000000001
$ws.Range("M120").Value = -62210.5
$ws.Range("N120").Value = -77054.20000000001
$ws.Range("H122").Value = 3968789.2
$ws.Range("I122").Value = 524.5714
$ws.Range("J122").Value = 9524360
$ws.Range("K122").Value = 4721.1426
$ws.Range("L122").Value = 85719240
$ws.Range("M122").Value = -2271.1426
$ws.Range("N122").Value = -85724140
$ws.Range("H129").Value = 19759202
$ws.Range("I129").Value = 1252.6364
$ws.Range("J129").Value = 33342792
$ws.Range("K129").Value = 3757.9092
$ws.Range("L129").Value = 100028376
$ws.Range("M129").Value = 1242.0908
$ws.Range("N129").Value = -100038376
$ws.Range("H131").Value = 14432528
$ws.Range("I131").Value = 10102587
$ws.Range("J131").Value = 22226422
$ws.Range("K131").Value = 30307761
$ws.Range("L131").Value = 66679266
$ws.Range("M131").Value = -30302721
$ws.Range("N131").Value = -66689346
$ws.Range("H134").Value = 2415.9473
$ws.Range("I134").Value = 1827.9445
$ws.Range("J134").Value = 13000
$ws.Range("K134").Value = 5483.833500000001
$ws.Range("L134").Value = 39000
$ws.Range("M134").Value = -413.8335000000006
$ws.Range("H137").Value = 1692.45
$ws.Range("I137").Value = 921.875
$ws.Range("K137").Value = 2765.625
$ws.Range("M137").Value = 2334.375
$ws.Range("H140").Value = 3701.2666
$ws.Range("I140").Value = 2804.182
$ws.Range("K140").Value = 8412.545999999998
$ws.Range("M140").Value = -3232.545999999998
$ws.Range("N82").Value = -57809
$ws.Range("N85").Value = -59805
$ws.Range("N134").Value = -49140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 129826.805
$ws.Range("I80").Value = 139572.03
$ws.Range("K80").Value = 139572.03
$ws.Range("M80").Value = -138574.03
$ws.Range("H83").Value = 129826.805
$ws.Range("I83").Value = 139572.03
$ws.Range("K83").Value = 697860.15
$ws.Range("M83").Value = -692868.15
$ws.Range("H107").Value = 54785.58
$ws.Range("I107").Value = 201255.8
$ws.Range("K107").Value = 201255.8
$ws.Range("M107").Value = -199335.8
$ws.Range("H113").Value = 4141.2856
$ws.Range("I113").Value = 3499.25
$ws.Range("J113").Value = 4997.3335
$ws.Range("K113").Value = 3499.25
$ws.Range("L113").Value = 4997.3335
$ws.Range("M113").Value = -1329.25
$ws.Range("N113").Value = -9337.333500000001
$ws.Range("H122").Value = 5900.4
$ws.Range("I122").Value = 6750.625
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 20251.875
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -17801.875
$ws.Range("N122").Value = -12398.5
$ws.Range("H139").Value = 91907.28999999999
$ws.Range("J139").Value = 91907.28999999999
$ws.Range("L139").Value = 91907.28999999999
$ws.Range("N139").Value = -102187.29
$ws.Range("H140").Value = 85722.25
$ws.Range("J140").Value = 87095.8
$ws.Range("L140").Value = 87095.8
$ws.Range("N140").Value = -97455.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2938.1052
$ws.Range("I40").Value = 2552.125
$ws.Range("K40").Value = 2552.125
$ws.Range("M40").Value = -2416.125
$ws.Range("H68").Value = 1766.3
$ws.Range("I68").Value = 1211.5
$ws.Range("K68").Value = 1211.5
$ws.Range("M68").Value = -462.5
$ws.Range("H71").Value = 1766.3
$ws.Range("I71").Value = 1211.5
$ws.Range("K71").Value = 6057.5
$ws.Range("M71").Value = -2313.5
$ws.Range("H122").Value = 2878.15
$ws.Range("I122").Value = 2730.0571
$ws.Range("J122").Value = 3914.8
$ws.Range("K122").Value = 8190.1713
$ws.Range("L122").Value = 11744.4
$ws.Range("M122").Value = -5740.1713
$ws.Range("N122").Value = -16644.4
$ws.Range("H136").Value = 2703.5
$ws.Range("I136").Value = 2350.4375
$ws.Range("K136").Value = 7051.3125
$ws.Range("M136").Value = -4501.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3573099.8
$ws.Range("J2").Value = 3573099.8
$ws.Range("L2").Value = 3573099.8
$ws.Range("N2").Value = -3573323.8
$ws.Range("H81").Value = 53181.65
$ws.Range("I81").Value = 2102.375
$ws.Range("J81").Value = 257498.75
$ws.Range("K81").Value = 4204.75
$ws.Range("L81").Value = 514997.5
$ws.Range("M81").Value = -3143.75
$ws.Range("N81").Value = -517119.5
$ws.Range("H84").Value = 53181.65
$ws.Range("I84").Value = 2102.375
$ws.Range("J84").Value = 257498.75
$ws.Range("K84").Value = 21023.75
$ws.Range("L84").Value = 2574987.5
$ws.Range("M84").Value = -15719.75
$ws.Range("N84").Value = -2585595.5
$ws.Range("H122").Value = 3357
$ws.Range("I122").Value = 2618.889
$ws.Range("K122").Value = 7856.667
$ws.Range("M122").Value = -5406.667
$ws.Range("H132").Value = 2347.762
$ws.Range("I132").Value = 1436.5385
$ws.Range("K132").Value = 4309.6155
$ws.Range("M132").Value = -1779.6155
